# Applies the Dec-09 Betfair odds refresh: updates back/lay odds and
# correct-score lay prices across rows 2-16 per the published diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.54
$ws.Range("G2").Value = 1.63
$ws.Range("I2").Value = 9.4
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 4.6
$ws.Range("L2").Value = 1.43
$ws.Range("N2").Value = 3.55
$ws.Range("O2").Value = 1.34
$ws.Range("P2").Value = 1.86
$ws.Range("Q2").Value = 2.04
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.7
$ws.Range("T2").Value = 2.04
$ws.Range("U2").Value = 1.81
$ws.Range("V2").Value = 1.13
$ws.Range("W2").Value = 2.6
$ws.Range("X2").Value = 1000

# Row 3
$ws.Range("F3").Value = 1.85
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 4.8
$ws.Range("I3").Value = 5.8
$ws.Range("J3").Value = 3.25
$ws.Range("K3").Value = 3.75
$ws.Range("L3").Value = 1.49
$ws.Range("N3").Value = 3.1
$ws.Range("O3").Value = 1.42
$ws.Range("Q3").Value = 2.26
$ws.Range("S3").Value = 4.3
$ws.Range("T3").Value = 1.99
$ws.Range("U3").Value = 1.84
$ws.Range("W3").Value = 2
$ws.Range("X3").Value = 11.5
$ws.Range("Y3").Value = 15.5
$ws.Range("Z3").Value = 1000
$ws.Range("AC3").Value = 8.4
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 11.5
$ws.Range("AG3").Value = 11
$ws.Range("AJ3").Value = 23
$ws.Range("AK3").Value = 24
$ws.Range("AL3").Value = 1000

# Row 4
$ws.Range("F4").Value = 1.79
$ws.Range("H4").Value = 5.9
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.7
$ws.Range("L4").Value = 1.56
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 2.78
$ws.Range("O4").Value = 1.48
$ws.Range("P4").Value = 1.58
$ws.Range("Q4").Value = 2.48
$ws.Range("R4").Value = 1.21
$ws.Range("S4").Value = 4.9
$ws.Range("T4").Value = 2.2
$ws.Range("U4").Value = 1.69
$ws.Range("V4").Value = 1.18
$ws.Range("W4").Value = 2.2
$ws.Range("X4").Value = 17
$ws.Range("AA4").Value = 260
$ws.Range("AB4").Value = 6.4
$ws.Range("AF4").Value = 21
$ws.Range("AG4").Value = 20
$ws.Range("AI4").Value = 170
$ws.Range("AK4").Value = 85
$ws.Range("AL4").Value = 370

# Row 5
$ws.Range("F5").Value = 7.2
$ws.Range("I5").Value = 1.55
$ws.Range("L5").Value = 1.37
$ws.Range("N5").Value = 4.5
$ws.Range("O5").Value = 1.26
$ws.Range("R5").Value = 1.47
$ws.Range("S5").Value = 3.05
$ws.Range("T5").Value = 1.94
$ws.Range("V5").Value = 2.82
$ws.Range("AC5").Value = 10
$ws.Range("AI5").Value = 32
$ws.Range("AK5").Value = 100
$ws.Range("AL5").Value = 90

# Row 6
$ws.Range("F6").Value = 1.24
$ws.Range("G6").Value = 1.25
$ws.Range("H6").Value = 14
$ws.Range("I6").Value = 14.5
$ws.Range("J6").Value = 8
$ws.Range("L6").Value = 1.22
$ws.Range("N6").Value = 8.4
$ws.Range("O6").Value = 1.12
$ws.Range("P6").Value = 3.45
$ws.Range("Q6").Value = 1.39
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 1.98
$ws.Range("T6").Value = 1.84
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.07
$ws.Range("W6").Value = 5
$ws.Range("Z6").Value = 150
$ws.Range("AA6").Value = 640
$ws.Range("AB6").Value = 14.5
$ws.Range("AC6").Value = 18
$ws.Range("AD6").Value = 48
$ws.Range("AE6").Value = 210
$ws.Range("AF6").Value = 9.800000000000001
$ws.Range("AG6").Value = 11.5
$ws.Range("AH6").Value = 29
$ws.Range("AI6").Value = 130
$ws.Range("AJ6").Value = 10.5
$ws.Range("AL6").Value = 28
$ws.Range("AM6").Value = 120
$ws.Range("AN6").Value = 3.25
$ws.Range("AO6").Value = 160

# Row 7
$ws.Range("F7").Value = 22
$ws.Range("G7").Value = 36
$ws.Range("H7").Value = 1.15
$ws.Range("I7").Value = 1.18
$ws.Range("J7").Value = 9
$ws.Range("K7").Value = 10.5
$ws.Range("L7").Value = 1.25
$ws.Range("M7").Value = 1.02
$ws.Range("N7").Value = 6.2
$ws.Range("O7").Value = 1.16
$ws.Range("P7").Value = 2.8
$ws.Range("Q7").Value = 1.47
$ws.Range("R7").Value = 1.73
$ws.Range("S7").Value = 2.2
$ws.Range("V7").Value = 6.6
$ws.Range("Y7").Value = 12.5
$ws.Range("Z7").Value = 9
$ws.Range("AA7").Value = 9.199999999999999
$ws.Range("AH7").Value = 65
$ws.Range("AI7").Value = 1000
$ws.Range("AO7").Value = 3.3

# Row 8
$ws.Range("F8").Value = 1.87
$ws.Range("G8").Value = 2.1
$ws.Range("H8").Value = 2.58
$ws.Range("L8").Value = 1.34
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 4.4
$ws.Range("O8").Value = 1.23
$ws.Range("P8").Value = 2.2
$ws.Range("Q8").Value = 1.68
$ws.Range("R8").Value = 1.47
$ws.Range("S8").Value = 2.66
$ws.Range("T8").Value = 1.64
$ws.Range("V8").Value = 1.04
$ws.Range("W8").Value = 1.9

# Row 9
$ws.Range("G9").Value = 3.45
$ws.Range("I9").Value = 2.32
$ws.Range("J9").Value = 3.65
$ws.Range("K9").Value = 3.7
$ws.Range("L9").Value = 1.35
$ws.Range("Q9").Value = 1.75
$ws.Range("S9").Value = 2.88
$ws.Range("V9").Value = 1.76
$ws.Range("W9").Value = 1.4
$ws.Range("X9").Value = 17.5
$ws.Range("AC9").Value = 8.199999999999999
$ws.Range("AF9").Value = 25
$ws.Range("AN9").Value = 26

# Row 10
$ws.Range("K10").Value = 5.2
$ws.Range("L10").Value = 1.37
$ws.Range("Q10").Value = 1.82
$ws.Range("S10").Value = 3.1
$ws.Range("T10").Value = 2.06
$ws.Range("X10").Value = 18.5
$ws.Range("AA10").Value = 280
$ws.Range("AC10").Value = 11
$ws.Range("AH10").Value = 26
$ws.Range("AM10").Value = 140

# Row 11
$ws.Range("J11").Value = 3.85
$ws.Range("P11").Value = 2.64
$ws.Range("S11").Value = 2.42
$ws.Range("T11").Value = 1.52
$ws.Range("V11").Value = 1.73

# Row 12
$ws.Range("F12").Value = 2.28
$ws.Range("G12").Value = 2.3
$ws.Range("I12").Value = 3.3
$ws.Range("J12").Value = 3.9
$ws.Range("L12").Value = 1.29
$ws.Range("Q12").Value = 1.56
$ws.Range("R12").Value = 1.7
$ws.Range("S12").Value = 2.38
$ws.Range("V12").Value = 1.43
$ws.Range("X12").Value = 24
$ws.Range("Y12").Value = 20
$ws.Range("AB12").Value = 16
$ws.Range("AJ12").Value = 30
$ws.Range("AL12").Value = 25

# Row 13
$ws.Range("F13").Value = 2.08
$ws.Range("G13").Value = 2.1
$ws.Range("H13").Value = 3.75
$ws.Range("I13").Value = 3.85
$ws.Range("J13").Value = 3.8
$ws.Range("K13").Value = 3.85
$ws.Range("L13").Value = 1.34
$ws.Range("O13").Value = 1.23
$ws.Range("P13").Value = 2.34
$ws.Range("Q13").Value = 1.72
$ws.Range("S13").Value = 2.8
$ws.Range("T13").Value = 1.63
$ws.Range("U13").Value = 2.52
$ws.Range("V13").Value = 1.35
$ws.Range("W13").Value = 1.89
$ws.Range("Y13").Value = 17
$ws.Range("Z13").Value = 29
$ws.Range("AA13").Value = 70
$ws.Range("AD13").Value = 15.5
$ws.Range("AE13").Value = 40
$ws.Range("AH13").Value = 15
$ws.Range("AJ13").Value = 25
$ws.Range("AO13").Value = 30

# Row 14
$ws.Range("H14").Value = 18
$ws.Range("I14").Value = 19
$ws.Range("J14").Value = 9.199999999999999
$ws.Range("K14").Value = 9.6
$ws.Range("L14").Value = 1.15
$ws.Range("N14").Value = 12
$ws.Range("P14").Value = 4.7
$ws.Range("Q14").Value = 1.25
$ws.Range("R14").Value = 2.5
$ws.Range("S14").Value = 1.63
$ws.Range("T14").Value = 1.73
$ws.Range("U14").Value = 2.28
$ws.Range("Y14").Value = 100
$ws.Range("AB14").Value = 20
$ws.Range("AC14").Value = 23
$ws.Range("AF14").Value = 12.5
$ws.Range("AG14").Value = 13
$ws.Range("AI14").Value = 140
$ws.Range("AM14").Value = 100
$ws.Range("AN14").Value = 2.44

# Row 15
$ws.Range("I15").Value = 2.46
$ws.Range("L15").Value = 1.35
$ws.Range("O15").Value = 1.25
$ws.Range("P15").Value = 2.26
$ws.Range("Q15").Value = 1.75
$ws.Range("R15").Value = 1.51
$ws.Range("S15").Value = 2.92
$ws.Range("T15").Value = 1.64
$ws.Range("V15").Value = 1.68
$ws.Range("Z15").Value = 16.5
$ws.Range("AD15").Value = 11
$ws.Range("AH15").Value = 15.5
$ws.Range("AM15").Value = 70

# Row 16
$ws.Range("F16").Value = 2.32
$ws.Range("G16").Value = 2.58
$ws.Range("H16").Value = 2.86
$ws.Range("I16").Value = 3.3
$ws.Range("J16").Value = 3.45
$ws.Range("L16").Value = 1.3
$ws.Range("N16").Value = 5.1
$ws.Range("O16").Value = 1.2
$ws.Range("P16").Value = 2.42
$ws.Range("Q16").Value = 1.64
$ws.Range("R16").Value = 1.59
$ws.Range("S16").Value = 2.58
$ws.Range("T16").Value = 1.5
$ws.Range("U16").Value = 2.56
$ws.Range("V16").Value = 1.45
$ws.Range("W16").Value = 1.64
$ws.Range("X16").Value = 21
$ws.Range("Y16").Value = 19
$ws.Range("Z16").Value = 24
$ws.Range("AA16").Value = 300
$ws.Range("AB16").Value = 17
$ws.Range("AC16").Value = 9
$ws.Range("AD16").Value = 15
$ws.Range("AE16").Value = 29
$ws.Range("AF16").Value = 19
$ws.Range("AK16").Value = 23
$ws.Range("AM16").Value = 60
$ws.Range("AN16").Value = 13
$ws.Range("AO16").Value = 20
